# Weekly update of Fruit/Vegetable (Alcachofa) prices in the "Mapocho Venta
# Directa de Santiago" sheet: each data row (2-14) is refreshed with a new
# reporting date (D), volume (J), min/max/weighted prices (K/L/M), origin
# (O) and $/Kg (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = fecha; J = volumen; K = precio minimo; L = precio maximo;
#           M = precio promedio ponderado; O = origen; P = precio $/Kg }
$rows = @{
    2  = @{ D = 44446; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 }
    3  = @{ D = 44421; J = 25; K = 15000; L = 16000; M = 15400; O = "Provincia de Limarí";  P = 513 }
    4  = @{ D = 44474; J = 45; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";  P = 333 }
    5  = @{ D = 44449; J = 45; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 }
    6  = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 }
    7  = @{ D = 44432; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí";  P = 467 }
    8  = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 }
    9  = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí";  P = 467 }
    10 = @{ D = 44460; J = 45; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí";  P = 433 }
    11 = @{ D = 44453; J = 50; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 }
    12 = @{ D = 44418; J = 30; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí";  P = 500 }
    13 = @{ D = 44425; J = 35; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 }
    14 = @{ D = 44376; J = 25; K = 18000; L = 18000; M = 18000; O = "Provincia de Limarí";  P = 600 }
}

foreach ($row in $rows.Keys) {
    $data = $rows[$row]
    $ws.Range("D$row").Value = $data.D
    $ws.Range("J$row").Value = $data.J
    $ws.Range("K$row").Value = $data.K
    $ws.Range("L$row").Value = $data.L
    $ws.Range("M$row").Value = $data.M
    $ws.Range("O$row").Value = $data.O
    $ws.Range("P$row").Value = $data.P
}
